$wb = $excel.ActiveWorkbook

# --- designs sheet: add 3 new rows (Output, Output efficiency, Output price) ---
$ws = $wb.Worksheets.Item("designs")
$ws.Range("A12").Value = "Class 8 Diesel Tractor"
$ws.Range("B12").Value = "Reference"
$ws.Range("C12").Value = "Output"
$ws.Range("D12").Value = "VMT"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "mile/mile"
$ws.Range("G12").Value = "Placeholder varible since output is needed"

$ws.Range("A13").Value = "Class 8 Diesel Tractor"
$ws.Range("B13").Value = "Reference"
$ws.Range("C13").Value = "Output efficiency"
$ws.Range("D13").Value = "VMT"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "mile/mile"
$ws.Range("G13").Value = "Placeholder varible since output is needed"

$ws.Range("A14").Value = "Class 8 Diesel Tractor"
$ws.Range("B14").Value = "Reference"
$ws.Range("C14").Value = "Output price"
$ws.Range("D14").Value = "VMT"
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = "$/mi"
$ws.Range("G14").Value = "Placeholder varible since output is needed"

$ws.Range("H15").Select()

# --- investments sheet: move the small table from M1:P2 to A1:D2 ---
$ws5 = $wb.Worksheets.Item("investments")
$ws5.Range("A1").Value = $ws5.Range("M1").Value
$ws5.Range("B1").Value = $ws5.Range("N1").Value
$ws5.Range("C1").Value = $ws5.Range("O1").Value
$ws5.Range("D1").Value = $ws5.Range("P1").Value
$ws5.Range("A2").Value = $ws5.Range("M2").Value
$ws5.Range("B2").Value = $ws5.Range("N2").Value
$ws5.Range("C2").Value = $ws5.Range("O2").Value

$ws5.Range("M1:P2").Clear()

$ws5.Range("D4").Select()

# --- parameters sheet: autofit columns A and B, change selection ---
$ws2 = $wb.Worksheets.Item("parameters")
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Range("C24").Select()

# --- indices sheet: just change selection to range A10:E10 ---
$ws4 = $wb.Worksheets.Item("indices")
$ws4.Range("A10:E10").Select()
